$wb = $excel.ActiveWorkbook

# Create a brand-new worksheet. A freshly created sheet never has any
# explicit <cols> column-width overrides recorded against it (unlike the
# existing "Sheet1", whose column widths can be resized but never fully
# cleared once set). Building the refreshed table on a clean sheet and
# then swapping it in for the old one is the only reliable way to drop
# the legacy custom column widths while keeping everything else intact.
$newWs = $wb.Worksheets.Add()
$oldWs = $wb.Worksheets.Item(2)

# Recreate the worksheet-level formatting the original sheet had.
$newWs.Outline.SummaryRow = 1
$newWs.Outline.SummaryColumn = 1

$ps = $newWs.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Header row: renamed columns, same bold / bordered / centered style as before.
$newWs.Range("A1").Value = "Date"
$newWs.Range("B1").Value = "Name"
$newWs.Range("C1").Value = "Entry Time"
$newWs.Range("D1").Value = "Exit Time"

$header = $newWs.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data row: updated entry - keep "2025-01-27" / time values as plain text,
# matching the original (not auto-converted to a date/number), and clear
# any residual number-format style picked up while forcing text entry.
$newWs.Range("A2").NumberFormat = "@"
$newWs.Range("A2").Value = "2025-01-27"
$newWs.Range("A2").Style = "Normal"

$newWs.Range("B2").Value = "Meharjot"

$newWs.Range("C2").NumberFormat = "@"
$newWs.Range("C2").Value = "18:47:06"
$newWs.Range("C2").Style = "Normal"

# Exit Time for the new entry is not known yet - leave it blank.
$newWs.Range("D2").Value = ""

# Drop the old sheet (with its legacy column widths) and put the new one
# in its place under the original name.
$oldWs.Delete() | Out-Null
$newWs.Name = "Sheet1"
